# DevTask@QRCode.xlsx - UI integration development plan
# Adds 10 new task rows (19-28, sheet rows 20-29) to Sheet1 and updates the
# "Status" column (H) for several existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# 1. Update the Status column (H) for existing rows 14-19
# ---------------------------------------------------------------------------
$ws.Range("H14").Value = "In Process"
$ws.Range("H15").Value = "Completed"
$ws.Range("H16").Value = "Completed"
$ws.Range("H17").Value = "Completed"
$ws.Range("H18").Value = "In Process"
$ws.Range("H19").Value = "Completed"

# ---------------------------------------------------------------------------
# 2. Add the new rows describing the UI integration development plan
# ---------------------------------------------------------------------------

# --- Row 20 (S.No 19) ---
$ws.Range("A20").Value = 19
$ws.Range("B20").Value = 'New Template creation : Impex for template, Content Slot Names, Content Slot for template (for both desktop and mobile)'
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 41963
$ws.Range("E20").Value = 41963
$ws.Range("F20").Value = 'Swapnil'

# --- Row 21 (S.No 20) ---
$ws.Range("A21").Value = 20
$ws.Range("B21").Value = 'Create Template Structure VM file'

# --- Row 22 (S.No 21) ---
$ws.Range("A22").Value = 21
$ws.Range("B22").Value = 'Create Template JSP, with HTML code integration(for both desktop and mobile)'

# --- Row 23 (S.No 22) ---
$ws.Range("A23").Value = 22
$ws.Range("B23").Value = 'Create CSROrderListComponent, CSROrderDetailsComponent in items.xml extending SimpleCMSComponent.'
$ws.Range("C23").Value = 0.5
$ws.Range("D23").Value = 41963
$ws.Range("E23").Value = 41963
$ws.Range("F23").Value = 'Swarnima'
$ws.Range("G23").Value = 'With System Update'

# --- Row 24 (S.No 23) ---
$ws.Range("A24").Value = 23
$ws.Range("B24").Value = 'Create impexes for page, content slots, content slots for page, components'
$ws.Range("C24").Value = 0.5
$ws.Range("D24").Value = 41963
$ws.Range("E24").Value = 41963
$ws.Range("F24").Value = 'Swarnima'

# --- Row 25 (S.No 24) ---
$ws.Range("A25").Value = 24
$ws.Range("B25").Value = 'Modify header.tag/jsp with the current UI(when HTML code is available)(for both desktop and mobile)'
$ws.Range("C25").Value = 0.5
$ws.Range("D25").Value = 41963
$ws.Range("E25").Value = 41963
$ws.Range("F25").Value = 'Prasun'

# --- Row 26 (S.No 25) ---
$ws.Range("A26").Value = 25
$ws.Range("B26").Value = 'Create CSR Orders Page Controller with JSP(integration with UI)(for both desktop and mobile)'
$ws.Range("C26").Value = 1
$ws.Range("D26").Value = 41964
$ws.Range("E26").Value = 41964
$ws.Range("F26").Value = 'Swarnima'
$ws.Range("G26").Value = 'Use code from existing PickInStoreOrder controller.'

# --- Row 27 (S.No 26) ---
$ws.Range("A27").Value = 26
$ws.Range("B27").Value = "Create Controller for CSR Orders List Component`nUI integration with CSR Orders List Page JSP(for both desktop and mobile)`n"
$ws.Range("C27").Value = 1.5
$ws.Range("D27").Value = 41963
$ws.Range("E27").Value = 41964
$ws.Range("F27").Value = 'Prasun'

# --- Row 28 (S.No 27) ---
$ws.Range("A28").Value = 27
$ws.Range("B28").Value = "Create CSR Order Details Component functionality for retrieving order and customer details through AJAX`nUI integration with CSR Order details Page JSP(when HTML code is available)(for both desktop and mobile)"
$ws.Range("C28").Value = 1.5
$ws.Range("D28").Value = 41964
$ws.Range("E28").Value = 41967
$ws.Range("F28").Value = 'Swarnima'

# --- Row 29 (S.No 28) ---
$ws.Range("A29").Value = 28
$ws.Range("B29").Value = 'Auto import of impexes during initialization'
$ws.Range("C29").Value = 1.5
$ws.Range("D29").Value = 41964
$ws.Range("E29").Value = 41967
$ws.Range("F29").Value = 'Swapnil'
$ws.Range("G29").Value = 'Add importing impexes methods as well as synchronize methods'

Write-Host "Values written"
